# Daily attendance processing - 2025-11-01 01:24:12
# Normalize specific "Recorded By" (column G) values on the attendance log:
# the automated "System" tag is moved from the front of the recorder list
# to the back for a known set of recorder-string patterns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, backup@backdoor.com, system" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Value2

    if ($null -eq $text) { continue }
    if (-not $map.ContainsKey($text)) { continue }

    $cell.Value2 = $map[$text]
}
